$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.9
$ws.Range("I3").Value = 2.75
$ws.Range("J3").Value = 3.75
$ws.Range("L3").Value = 3.5
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 2.2
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 13
$ws.Range("AE3").Value = 17
$ws.Range("AG3").Value = 6.5
$ws.Range("AN3").Value = 4.75
$ws.Range("AO3").Value = 19
$ws.Range("AQ3").Value = 67
$ws.Range("AS3").Value = 351
$ws.Range("AT3").Value = 2.2
$ws.Range("AU3").Value = 9.5
$ws.Range("G4").Value = 1.75
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 5.5
$ws.Range("J4").Value = 2.5
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("U4").Value = 2.5
$ws.Range("V4").Value = 1.5
$ws.Range("W4").Value = 4.75
$ws.Range("X4").Value = 6.5
$ws.Range("Y4").Value = 9.5
$ws.Range("Z4").Value = 13
$ws.Range("AA4").Value = 19
$ws.Range("AC4").Value = 6
$ws.Range("AG4").Value = 10
$ws.Range("AH4").Value = 26
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 67
$ws.Range("AK4").Value = 51
$ws.Range("AL4").Value = 67
$ws.Range("AN4").Value = 3.5
$ws.Range("AO4").Value = 10
$ws.Range("AR4").Value = 67
$ws.Range("AW4").Value = 6.5
$ws.Range("AX4").Value = 34
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 201
$ws.Range("H7").Value = 2.9
$ws.Range("J7").Value = 2.85
$ws.Range("K7").Value = 1.98
$ws.Range("N7").Value = 6.7
$ws.Range("T7").Value = 2.42
$ws.Range("U7").Value = 1.83
$ws.Range("W7").Value = 6.4
$ws.Range("X7").Value = 10.25
$ws.Range("AC7").Value = 7.2
$ws.Range("AD7").Value = 5.7
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 8.75
$ws.Range("AH7").Value = 17.5
$ws.Range("AP7").Value = 20
$ws.Range("AS7").Value = 250
$ws.Range("AT7").Value = 2.4
$ws.Range("G14").Value = 3.2
$ws.Range("H14").Value = 3.1
$ws.Range("I14").Value = 2.2
$ws.Range("J14").Value = 3.7
$ws.Range("K14").Value = 2.05
$ws.Range("L14").Value = 2.8
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 9.37
$ws.Range("O14").Value = 1.35
$ws.Range("P14").Value = 2.7
$ws.Range("Q14").Value = 2.02
$ws.Range("R14").Value = 1.62
$ws.Range("S14").Value = 1.4
$ws.Range("T14").Value = 2.5
$ws.Range("U14").Value = 1.8
$ws.Range("V14").Value = 1.82
$ws.Range("W14").Value = 9
$ws.Range("X14").Value = 16.5
$ws.Range("Y14").Value = 11.25
$ws.Range("Z14").Value = 45
$ws.Range("AA14").Value = 30
$ws.Range("AB14").Value = 40
$ws.Range("AC14").Value = 8.25
$ws.Range("AD14").Value = 6
$ws.Range("AE14").Value = 14.5
$ws.Range("AF14").Value = 75
$ws.Range("AG14").Value = 6.8
$ws.Range("AH14").Value = 10
$ws.Range("AI14").Value = 9
$ws.Range("AJ14").Value = 21
$ws.Range("AK14").Value = 19.5
$ws.Range("AL14").Value = 32
$ws.Range("AN14").Value = 5
$ws.Range("AO14").Value = 17.5
$ws.Range("AP14").Value = 24
$ws.Range("AQ14").Value = 90
$ws.Range("AR14").Value = 120
$ws.Range("AS14").Value = 300
$ws.Range("AT14").Value = 2.47
$ws.Range("AU14").Value = 6.9
$ws.Range("AV14").Value = 60
$ws.Range("AW14").Value = 4.05
$ws.Range("AX14").Value = 11.5
$ws.Range("AY14").Value = 20
$ws.Range("AZ14").Value = 45
$ws.Range("BA14").Value = 80
